# Add a new "grupo" column (G) to the "Total" sheet, with the chemical
# group number for each element, as described in the commit message:
# "...añadir a la hoja de cálculo el grupo de cada elemento para
# identificar los ácidos"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(1, 7).Value = "grupo"

# Group values for rows 2..45 (element rows), column G
$grupos = @(1,1,15,11,13,2,2,17,14,20,12,17,9,6,1,11,17,8,1,13,1,12,17,1,1,2,7,15,1,10,16,15,14,10,10,2,1,16,16,14,14,2,4,12)

for ($i = 0; $i -lt $grupos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $grupos[$i]
}

# Update the frozen-pane / selection view state to match the target
$ws.Application.ActiveWindow.SplitColumn = 0
$ws.Application.ActiveWindow.SplitRow = 1

$ws.Range("G45").Select()
